$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 16: swap the table style (tableStyleId) applied to the
#    2-column cash-flow table.
# ---------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
$tableShape = $slide16.Shapes.Item(3)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{B7D272B4-8AA1-4471-8289-07790F5FD95D}")
}

# ---------------------------------------------------------------------------
# 2) Theme colours: the deck's primary theme (theme1.xml, currently the
#    "Integral" palette) is swapped for the default "Office Theme" palette.
# ---------------------------------------------------------------------------
function Set-ThemeRGB($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme

Set-ThemeRGB $themeColors 1 "000000"
Set-ThemeRGB $themeColors 2 "FFFFFF"
Set-ThemeRGB $themeColors 3 "44546A"
Set-ThemeRGB $themeColors 4 "E7E6E6"
Set-ThemeRGB $themeColors 5 "5B9BD5"
Set-ThemeRGB $themeColors 6 "ED7D31"
Set-ThemeRGB $themeColors 7 "A5A5A5"
Set-ThemeRGB $themeColors 8 "FFC000"
Set-ThemeRGB $themeColors 9 "4472C4"
Set-ThemeRGB $themeColors 10 "70AD47"
Set-ThemeRGB $themeColors 11 "0563C1"
Set-ThemeRGB $themeColors 12 "954F72"


